# Insert a new "idx" worksheet right before "sumup" and populate it with
# index/FX/commodity ticker rows, mirroring how the data was originally
# typed in (a few rows by hand, then the "group" column filled down, then
# the remaining rows pasted in).

$wb = $excel.ActiveWorkbook
$sumup = $wb.Worksheets.Item("sumup")
$ws = $wb.Worksheets.Add($sumup)
$ws.Name = "idx"

# --- Phase 1: type the ticker (column A) for the first three rows ---
$ws.Cells.Item(1, 1).Value = "^FCHI"
$ws.Cells.Item(2, 1).Value = "^GSPC"
$ws.Cells.Item(3, 1).Value = "^DJI"

# --- Phase 2: fill the "group" column (D) down for all 25 rows ---
for ($r = 1; $r -le 25; $r++) {
    $ws.Cells.Item($r, 4).Value = "IDX"
}

# --- Phase 3: type the country/market (column B) for the first three rows ---
$ws.Cells.Item(1, 2).Value = "FR"
$ws.Cells.Item(2, 2).Value = "US"
$ws.Cells.Item(3, 2).Value = "US"

# --- Phase 4: remaining rows 4-25, ticker then country ---
$tickers = @("^IXIC", "^FTSE", "^GDAXI", "^N225", "CLM15.NYM", "^BVSP", "^GSPTSE", "^MERV", "^MXX", "^ATX", "^BFX", "FTSEMIB.MI", "^SSMI", "^HIS", "^SBF250", "EURUSD=X", "EURBRL=X", "GBP=X", "EURCHF=X", "XAUUSD=X", "HGJ15.CMX", "^XAU")
$countries = @("US", "UK", "GE", "JP", "US", "BR", "CN", "US", "US", "FR", "FR", "FR", "US", "US", "FR", "FR", "FR", "FR", "FR", "US", "US", "US")
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $r = $i + 4
    $ws.Cells.Item($r, 1).Value = $tickers[$i]
    $ws.Cells.Item($r, 2).Value = $countries[$i]
}

# --- Phase 5: "active" flag (column E) for all 25 rows ---
for ($r = 1; $r -le 25; $r++) {
    $ws.Cells.Item($r, 5).Value = $true
}

# Leave the selection the way it ended up in the authored workbook, and make
# sure the new sheet is the active / visible tab.
$ws.Activate()
$ws.Range("C27").Select()
